$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.423895359039307
$ws.Range("B1").Value = 2.530923128128052
$ws.Range("C1").Value = 2.908759593963623
$ws.Range("D1").Value = 4.602647304534912
$ws.Range("E1").Value = 4.5048508644104
